$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 45.068240000000003
$ws.Range("B3").Value = 65.918144666666649
$ws.Range("B4").Value = 24.894326000000003
$ws.Range("B5").Value = 79.331043999999991
$ws.Range("B6").Value = 91.721729333333329
$ws.Range("B7").Value = 52.801807999999994
$ws.Range("B8").Value = 44.974846000000007
$ws.Range("B9").Value = 8.5249053333333347
